$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume updates. Values in column D (Price) that look like plain
# decimals (e.g. "1.001") would otherwise be auto-converted to numbers by the
# COM layer, losing the original text formatting (trailing zeros, etc.), so
# those cells are explicitly forced to Text before the value is written.

$ws.Range("D2").Value = "24.593.20"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "1.692.76"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.12"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3944"
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4014"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.520"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.000"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.53"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08763"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.211"
$ws.Range("E13").Value = "  +5.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.23"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.188"
$ws.Range("E15").Value = "  +11.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001312"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "1.697.68"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.74"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07065"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.65"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.061"
$ws.Range("E21").Value = "  +6.05%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.20"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").Value = "24.606.38"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.127"
$ws.Range("E25").Value = "  +8.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.335"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.80"
$ws.Range("E27").Value = "  +4.42%  "
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.45"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.196"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.469"
$ws.Range("E31").Value = "  +9.10%  "
$ws.Range("D32").Value = "1.884.25"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.076"
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08599"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.182"
$ws.Range("E35").Value = "  +6.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.51"
$ws.Range("E36").Value = "  +9.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2734"
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.925"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.39"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09113"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02731"
$ws.Range("E41").Value = "  +7.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.489"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7649"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.595"
$ws.Range("E44").Value = "  +6.83%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7150"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.227"
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.77"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.320"
$ws.Range("E50").Value = "  +7.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07985"
$ws.Range("E51").Value = "  +1.52%  "
